$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new task rows into the "Registration"/Tracker section (rows 84-88) ---
# Original layout:
#   84: * Python Tracker Connection   C=4
#   85: * Tracker GUI                 C=4
#   86: * Tracker DataCapture         C=3   -> becomes 5
#   87: * Registration                C=3
#   88: Total Hours  =SUM(C84:C87)=14
#   89: Paid    D=4
#   90: Not Paid D=10 -> becomes 15
#
# Target layout:
#   84: * Python Tracker Connection   C=4
#   85: * Tracker GUI                 C=4
#   86: * Tracker DataCapture         C=5
#   87: * 2D/3D Views     (NEW)       C=2
#   88: * Registration                C=3
#   89: * Video Capture   (NEW)       C=1
#   90: Total Hours  =SUM(C84:C89)=19
#   91: Paid    D=4
#   92: Not Paid D=15

# Insert a blank row at 87 (pushes "Registration" 87->88 and the Total row 88->89)
$ws.Rows("87:87").Insert()
$ws.Range("B86:C86").Copy()
$ws.Range("B87:C87").PasteSpecial(-4122)

# Insert a blank row at 89 (pushes the Total row 89->90), leaving "Registration" at 88
$ws.Rows("89:89").Insert()
$ws.Range("B88:C88").Copy()
$ws.Range("B89:C89").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new shared-string-backed values. "* Video Capture" must be
# introduced into the shared-string table before "* 2D/3D Views" (matches
# the original authoring order), so set that cell's text first.
$ws.Cells.Item(89, 2).Value = "* Video Capture"
$ws.Cells.Item(87, 2).Value = "* 2D/3D Views"

$ws.Cells.Item(87, 3).Value = 2
$ws.Cells.Item(89, 3).Value = 1

# Update existing values that changed
$ws.Cells.Item(86, 3).Value = 5

# Fix up the Total Hours formula (now at row 90) to cover the new range
$ws.Cells.Item(90, 3).Formula = "=SUM(C84:C89)"

# Not Paid hours total (row shifted from 90 -> 92)
$ws.Cells.Item(92, 4).Value = 15

# Match the final cell selection recorded in the saved file
[void]$ws.Range("M96:R96").Select()

$wb.Save()
